$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''51.518.91'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.99%  '
$ws.Range('D3').Value = '''2.982.94'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.42%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''382.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.88%  '
$ws.Range('D6').Value = '''104.23'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.25%  '
$ws.Range('D7').Value = '''0.546'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.01%  '
$ws.Range('D8').Value = '''0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '''0.595'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.18%  '
$ws.Range('D10').Value = '''37.31'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.21%  '
$ws.Range('E11').Value = '  +0.13%  '
$ws.Range('D12').Value = '''0.0846'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.46%  '
$ws.Range('D13').Value = '''3.450.51'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.46%  '
$ws.Range('D14').Value = '''18.41'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.73%  '
$ws.Range('D15').Value = '''7.57'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.01%  '
$ws.Range('D16').Value = '''2.986.21'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.55%  '
$ws.Range('D17').Value = '''0.974'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.76%  '
$ws.Range('D18').Value = '''51.446.98'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.00%  '
$ws.Range('E19').Value = '  +2.78%  '
$ws.Range('D20').Value = '''7.45'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.01%  '
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('D22').Value = '''0.0₃0965'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.54%  '
$ws.Range('D23').Value = '''68.94'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.44%  '
$ws.Range('D24').Value = '''262.90'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.72%  '
$ws.Range('D25').Value = '''2.92'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +9.08%  '
$ws.Range('D26').Value = '''8.34'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +18.07%  '
$ws.Range('D27').Value = '''7.74'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +23.53%  '
$ws.Range('E28').Value = '  +14.01%  '
$ws.Range('D29').Value = '''0.170'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.11%  '
$ws.Range('D30').Value = '''26.01'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.81%  '
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('E32').Value = '  +0.60%  '
$ws.Range('D33').Value = '''34.83'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Value = '''51.01'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.51%  '
$ws.Range('D35').Value = '''2.07'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.05%  '
$ws.Range('E36').Value = '  +7.52%  '
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').Value = '''3.03'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.73%  '
$ws.Range('D39').Value = '''17.08'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.56%  '
$ws.Range('E40').Value = '  -0.08%  '
$ws.Range('E41').Value = '  +0.43%  '
$ws.Range('E42').Value = '  +3.16%  '
$ws.Range('D43').Value = '''122.11'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.33%  '
$ws.Range('D44').Value = '''21.76'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.99%  '
$ws.Range('D45').Value = '''0.280'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +18.99%  '
$ws.Range('E46').Value = '  -2.64%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '''3.28'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.99%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '''2.032.50'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.59%  '
$ws.Range('D50').Value = '''0.0332'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.95%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').Value = '''58.18'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.35%  '
